# Cotações atualizadas - 2025-10-18
# Append a new daily quotation row (row 44) to the bottom of the table,
# carrying forward the previous day's fund values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A: new date, formatted the same way as the existing date column (A2:A43)
$dateCell = $ws.Range("A44")
$dateCell.Value = 45948
$dateCell.NumberFormat = $ws.Range("A43").NumberFormat

# Columns B-E: quotation values (carried forward from the prior day, 2025-10-17)
$ws.Range("B44").Value = "21,7414"
$ws.Range("C44").Value = "15,4996"
$ws.Range("D44").Value = "15,5859"
$ws.Range("E44").Value = "15,5859"
